# Auto-applied bulk value update ("chore: update Sheets via scheduled runner").
# Each sheet's Leve-profit table (columns H:N) is refreshed with newly-pulled
# market-board prices/profits. H:J are per-unit price lookups, K:N are the
# leve totals (scaled by "Leve Amount"); a handful of rows gain or lose a
# trailing N cell because the source feed omits blank/zero results.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 5393.1875
$ws.Range("I40").Value = 2796.6667
$ws.Range("J40").Value = 5992.385
$ws.Range("K40").Value = 2796.6667
$ws.Range("L40").Value = 5992.385
$ws.Range("M40").Value = -2621.6667
$ws.Range("H98").Value = 27270
$ws.Range("I98").Value = 28187.572
$ws.Range("J98").Value = 8001
$ws.Range("K98").Value = 28187.572
$ws.Range("L98").Value = 8001
$ws.Range("M98").Value = -26689.572
$ws.Range("N98").Value = -10997
$ws.Range("H122").Value = 27270
$ws.Range("I122").Value = 28187.572
$ws.Range("J122").Value = 8001
$ws.Range("K122").Value = 84562.716
$ws.Range("L122").Value = 24003
$ws.Range("M122").Value = -82112.716
$ws.Range("N122").Value = -28903
$ws.Range("H123").Value = 145555
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 145555
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 145555
$ws.Range("N123").Value = -155355
$ws.Range("H137").Value = 710975.75
$ws.Range("I137").Value = 1064673.8
$ws.Range("J137").Value = 3579.6667
$ws.Range("K137").Value = 3194021.4
$ws.Range("L137").Value = 10739.0001
$ws.Range("M137").Value = -3191471.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7666.7466
$ws.Range("I32").Value = 7666.7466
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 7666.7466
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -7379.7466

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 74486.3
$ws.Range("I26").Value = 60838.145
$ws.Range("J26").Value = 106332
$ws.Range("K26").Value = 60838.145
$ws.Range("L26").Value = 106332
$ws.Range("M26").Value = -60546.145
$ws.Range("H52").Value = 79750
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 79750
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 79750
$ws.Range("N52").Value = -80276
$ws.Range("H55").Value = 82250
$ws.Range("I55").Value = 90000
$ws.Range("J55").Value = 74500
$ws.Range("K55").Value = 90000
$ws.Range("L55").Value = 74500
$ws.Range("M55").Value = -89727
$ws.Range("N55").Value = -75046
$ws.Range("H61").Value = 44998.5
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 44998.5
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 44998.5
$ws.Range("N61").Value = -45624.5
$ws.Range("H99").Value = 24580.555
$ws.Range("I99").Value = 32555.166
$ws.Range("J99").Value = 8631.333000000001
$ws.Range("K99").Value = 32555.166
$ws.Range("L99").Value = 8631.333000000001
$ws.Range("M99").Value = -31057.166
$ws.Range("H105").Value = 10008.944
$ws.Range("I105").Value = 12971.9
$ws.Range("J105").Value = 6305.25
$ws.Range("K105").Value = 12971.9
$ws.Range("L105").Value = 6305.25
$ws.Range("M105").Value = -11224.9
$ws.Range("N105").Value = -9799.25
$ws.Range("H121").Value = 79750
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 79750
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 79750
$ws.Range("N121").Value = -83244

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2673.491
$ws.Range("I31").Value = 1834.1025
$ws.Range("J31").Value = 4719.5
$ws.Range("K31").Value = 1834.1025
$ws.Range("L31").Value = 4719.5
$ws.Range("M31").Value = -1539.1025
$ws.Range("N31").Value = -5309.5
$ws.Range("H34").Value = 2673.491
$ws.Range("I34").Value = 1834.1025
$ws.Range("J34").Value = 4719.5
$ws.Range("K34").Value = 1834.1025
$ws.Range("L34").Value = 4719.5
$ws.Range("M34").Value = -1632.1025
$ws.Range("N34").Value = -5123.5
$ws.Range("H54").Value = 66345.39999999999
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 66345.39999999999
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 66345.39999999999
$ws.Range("N54").Value = -67661.39999999999
$ws.Range("H99").Value = 139141.19
$ws.Range("I99").Value = 220489.39
$ws.Range("J99").Value = 5497.7144
$ws.Range("K99").Value = 220489.39
$ws.Range("L99").Value = 5497.7144
$ws.Range("M99").Value = -218991.39
$ws.Range("H105").Value = 4340.811
$ws.Range("I105").Value = 7678.684
$ws.Range("J105").Value = 817.5
$ws.Range("K105").Value = 7678.684
$ws.Range("L105").Value = 817.5
$ws.Range("M105").Value = -5931.684
$ws.Range("N105").Value = -4311.5
$ws.Range("H109").Value = 69999.5
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 69999.5
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 69999.5
$ws.Range("N109").Value = -72079.5
$ws.Range("H114").Value = 57554.668
$ws.Range("I114").Value = 0
$ws.Range("J114").Value = 57554.668
$ws.Range("K114").Value = 0
$ws.Range("L114").Value = 57554.668
$ws.Range("N114").Value = -66232.66800000001
$ws.Range("H126").Value = 139141.19
$ws.Range("I126").Value = 220489.39
$ws.Range("J126").Value = 5497.7144
$ws.Range("K126").Value = 661468.17
$ws.Range("L126").Value = 16493.1432
$ws.Range("M126").Value = -658998.17
$ws.Range("H132").Value = 3840.6572
$ws.Range("I132").Value = 3888.9119
$ws.Range("J132").Value = 2200
$ws.Range("K132").Value = 11666.7357
$ws.Range("L132").Value = 6600
$ws.Range("M132").Value = -9136.735700000001
$ws.Range("N132").Value = -11660
$ws.Range("H141").Value = 326107.44
$ws.Range("I141").Value = 56912.832
$ws.Range("J141").Value = 393406.1
$ws.Range("K141").Value = 56912.832
$ws.Range("L141").Value = 393406.1
$ws.Range("M141").Value = -51732.832
$ws.Range("N141").Value = -403766.1

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H127").Value = 1516.5
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 1516.5
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 4549.5
$ws.Range("N127").Value = -14469.5
$ws.Range("H140").Value = 3335075
$ws.Range("I140").Value = 3335075
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 10005225
$ws.Range("L140").Value = 0
$ws.Range("M140").Value = -10000045
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 28166.666
$ws.Range("I46").Value = 7000
$ws.Range("J46").Value = 49333.332
$ws.Range("K46").Value = 7000
$ws.Range("L46").Value = 49333.332
$ws.Range("M46").Value = -6844
$ws.Range("N46").Value = -49645.332

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 33241.332
$ws.Range("I40").Value = 72891.164
$ws.Range("J40").Value = 13416.417
$ws.Range("K40").Value = 72891.164
$ws.Range("L40").Value = 13416.417
$ws.Range("M40").Value = -72755.164
$ws.Range("N40").Value = -13688.417
$ws.Range("H46").Value = 2846.3845
$ws.Range("I46").Value = 800
$ws.Range("J46").Value = 5233.8335
$ws.Range("K46").Value = 800
$ws.Range("L46").Value = 5233.8335
$ws.Range("M46").Value = -612
$ws.Range("N46").Value = -5609.8335
$ws.Range("H55").Value = 5247.222
$ws.Range("I55").Value = 1824.5
$ws.Range("J55").Value = 6225.143
$ws.Range("K55").Value = 1824.5
$ws.Range("L55").Value = 6225.143
$ws.Range("M55").Value = -1651.5
$ws.Range("N55").Value = -6571.143
$ws.Range("H92").Value = 77500
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 77500
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 77500
$ws.Range("N92").Value = -82492
$ws.Range("H122").Value = 6150.385
$ws.Range("I122").Value = 5806.875
$ws.Range("J122").Value = 6700
$ws.Range("K122").Value = 17420.625
$ws.Range("L122").Value = 20100
$ws.Range("M122").Value = -14970.625
$ws.Range("H123").Value = 107493
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 107493
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 107493
$ws.Range("N123").Value = -117293

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 20281.688
$ws.Range("I107").Value = 1745.909
$ws.Range("J107").Value = 61060.4
$ws.Range("K107").Value = 5237.727000000001
$ws.Range("L107").Value = 183181.2
$ws.Range("M107").Value = -3317.727000000001
$ws.Range("N107").Value = -187021.2
$ws.Range("H122").Value = 5381.433
$ws.Range("I122").Value = 3371.5
$ws.Range("J122").Value = 7678.5
$ws.Range("K122").Value = 10114.5
$ws.Range("L122").Value = 23035.5
$ws.Range("M122").Value = -7664.5
$ws.Range("N122").Value = -27935.5
$ws.Range("H126").Value = 36049
$ws.Range("I126").Value = 46454.223
$ws.Range("J126").Value = 4833.3335
$ws.Range("K126").Value = 139362.669
$ws.Range("L126").Value = 14500.0005
$ws.Range("M126").Value = -136892.669
$ws.Range("N126").Value = -19440.0005
